$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.245.06"
$ws.Range("E2").Value = "  -3.71%  "
$ws.Range("D3").Value = "2.602.66"
$ws.Range("E3").Value = "  -0.96%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "'504.88"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.45%  "
$ws.Range("D6").Value = "'144.77"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.43%  "
$ws.Range("D7").Value = "'0.997"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("D8").Value = "'0.563"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.05%  "
$ws.Range("D9").Value = "2.629.68"
$ws.Range("E9").Value = "  -0.24%  "
$ws.Range("E10").Value = "  +0.29%  "
$ws.Range("E11").Value = "  -3.11%  "
$ws.Range("E12").Value = "  -3.18%  "
$ws.Range("D13").Value = "'0.126"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.45%  "
$ws.Range("D14").Value = "3.064.53"
$ws.Range("E14").Value = "  -0.77%  "
$ws.Range("D15").Value = "58.287.77"
$ws.Range("E15").Value = "  -3.68%  "
$ws.Range("D16").Value = "'20.94"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.67%  "
$ws.Range("E17").Value = "  -2.40%  "
$ws.Range("D18").Value = "2.619.88"
$ws.Range("E18").Value = "  -0.70%  "
$ws.Range("D19").Value = "'4.52"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Value = "'339.43"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.40%  "
$ws.Range("D21").Value = "'10.29"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.75%  "
$ws.Range("D22").Value = "'6.06"
$ws.Range("D22").Style = "Normal"
$ws.Range("E23").Value = "  +0.57%  "
$ws.Range("D24").Value = "'60.36"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.16%  "
$ws.Range("E25").Value = "  -0.74%  "
$ws.Range("D26").Value = "2.757.70"
$ws.Range("E26").Value = "  +0.67%  "
$ws.Range("D27").Value = "'0.999"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.01%  "
$ws.Range("D28").Value = "'0.158"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.70%  "
$ws.Range("D29").Value = "0.0₃0808"
$ws.Range("E29").Value = "  -3.24%  "
$ws.Range("D30").Value = "'6.99"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.98%  "
$ws.Range("E31").Value = "  -0.19%  "
$ws.Range("D32").Value = "'6.48"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +6.99%  "
$ws.Range("D33").Value = "'18.73"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.86%  "
$ws.Range("D34").Value = "'1.55"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.97%  "
$ws.Range("D35").Value = "'148.30"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.27%  "
$ws.Range("E36").Value = "  +14.42%  "
$ws.Range("D37").Value = "'3.96"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.29%  "
$ws.Range("E38").Value = "  -2.87%  "
$ws.Range("D39").Value = "'0.851"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.47%  "
$ws.Range("D40").Value = "'36.14"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.50%  "
$ws.Range("E41").Value = "  -1.11%  "
$ws.Range("D42").Value = "'1.39"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.29%  "
$ws.Range("D43").Value = "'0.615"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.05%  "
$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D44").Value = "'0.997"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.12%  "
$ws.Range("B45").Value = "Bittensor"
$ws.Range("C45").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D45").Value = "'276.15"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.02%  "
$ws.Range("D46").Value = "'0.0983"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.41%  "
$ws.Range("D47").Value = "'19.32"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.56%  "
$ws.Range("E48").Value = "  -3.35%  "
$ws.Range("B49").Value = "WhiteBITCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D49").Value = "'10.24"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.42%  "
$ws.Range("B50").Value = "VeChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D50").Value = "'0.0228"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.33%  "
$ws.Range("D51").Value = "'4.68"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.26%  "
